# Generate Report for Handoff
# Replaces the e2e test markdown files referenced by the localization status
# report with a new pair of files, and flips the rows from "handed back"
# state to a freshly generated "ready for handoff" state.

$wb = $excel.ActiveWorkbook

$oldFile1 = "27f5425e-7377-4959-9110-1f54699a9831"
$newFile1 = "721691ad-7b3d-415f-b730-8c9ee5d775ff"
$oldFile2 = "2f7db598-7b81-4391-a6ab-0a7ed8fae673"
$newFile2 = "ffff30423dfa-374b-4c3f-9829-48fe5f458970"

$newHash = "25efd87e78c42f047c0dc7ff0d780f1539a66c9d"

$newStatus = "Ready for handoff"
$newGenDate = "2016-08-25 09:07:35"
$newHandoffDateZh = "2016-08-25 09:07:30"
$newHandoffDateDe = "2016-08-25 09:07:35"
$nullDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newFile1.md"
$ws.Range("B2").Value = "e2e\$newFile1.md"
$ws.Range("E2").Value = $newStatus
$ws.Range("F2").Value = $newStatus
$ws.Range("G2").Value = $newGenDate

$ws.Range("A3").Value = "$newFile2.md"
$ws.Range("B3").Value = "e2e\$newFile2.md"
$ws.Range("E3").Value = $newStatus
$ws.Range("F3").Value = $newStatus
$ws.Range("G3").Value = $newGenDate

foreach ($hl in @($ws.Hyperlinks)) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newFile1.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\$newFile2.md"
    }
}

$ws.Columns.Item(5).ColumnWidth = 16.333333333333332
$ws.Columns.Item(6).ColumnWidth = 16.333333333333332

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$newFile1.md"
$ws.Range("C2").Value = $newStatus
$ws.Range("G2").Value = "$newFile1.$newHash.zh-cn.xlf"
$ws.Range("H2").Value = $newHandoffDateZh
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = $nullDate

$ws.Range("A3").Value = "$newFile2.md"
$ws.Range("C3").Value = $newStatus
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = "$newFile1.$newHash.zh-cn.xlf"
$ws.Range("H3").Value = $newHandoffDateZh
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = $nullDate

foreach ($hl in @($ws.Hyperlinks)) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$I$2' -or $addr -eq '$I$3') {
        $hl.Delete()
    }
}
foreach ($hl in @($ws.Hyperlinks)) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newFile1.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "$newFile2.md"
    }
}

$ws.Range("I2").Style = "Normal"
$ws.Range("I3").Style = "Normal"

$ws.Columns.Item(3).ColumnWidth = 16.333333333333332
$ws.Columns.Item(9).ColumnWidth = 17.833333333333332
$ws.Columns.Item(10).ColumnWidth = 20.833333333333332

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$newFile1.md"
$ws.Range("C2").Value = $newStatus
$ws.Range("G2").Value = "$newFile1.$newHash.de-de.xlf"
$ws.Range("H2").Value = $newHandoffDateDe
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = $nullDate

$ws.Range("A3").Value = "$newFile2.md"
$ws.Range("C3").Value = $newStatus
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = "$newFile1.$newHash.de-de.xlf"
$ws.Range("H3").Value = $newHandoffDateDe
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = $nullDate

foreach ($hl in @($ws.Hyperlinks)) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$I$2' -or $addr -eq '$I$3') {
        $hl.Delete()
    }
}
foreach ($hl in @($ws.Hyperlinks)) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newFile1.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "$newFile2.md"
    }
}

$ws.Range("I2").Style = "Normal"
$ws.Range("I3").Style = "Normal"

$ws.Columns.Item(3).ColumnWidth = 16.333333333333332
$ws.Columns.Item(9).ColumnWidth = 17.833333333333332
$ws.Columns.Item(10).ColumnWidth = 20.833333333333332
